$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "UserCreation" sheet: Sample Management Finished Product location/dept
#    code update.  Plant-1 -> Plant-3 for the automation test users, and the
#    related Department codes move from Plant-1-QC/Plant-1-QA to the new
#    Plant-3-QC / Plant-3-PF / Plant-3-QA codes.
# ---------------------------------------------------------------------------
$wsUser = $wb.Worksheets.Item("UserCreation")

# Location column (B) : Plant-1 -> Plant-3 for rows 2-9
$wsUser.Range("B2").Value = "Plant-3"
$wsUser.Range("B3").Value = "Plant-3"
$wsUser.Range("B4").Value = "Plant-3"
$wsUser.Range("B5").Value = "Plant-3"
$wsUser.Range("B6").Value = "Plant-3"
$wsUser.Range("B7").Value = "Plant-3"
$wsUser.Range("B8").Value = "Plant-3"
$wsUser.Range("B9").Value = "Plant-3"

# Department column (F)
$wsUser.Range("F2").Value = "Plant-3-QC"
$wsUser.Range("F3").Value = "Plant-3-QC"
$wsUser.Range("F4").Value = "Plant-3-PF"
$wsUser.Range("F5").Value = "Plant-3-PF"
$wsUser.Range("F6").Value = "Plant-3-QC"
$wsUser.Range("F7").Value = "Plant-3-QC"
$wsUser.Range("F8").Value = "Plant-3-QA"
$wsUser.Range("F9").Value = "Plant-3-QA"

# ---------------------------------------------------------------------------
# 2) "TestDetails" sheet: reorder so Qualitative tests are listed before the
#    Quantitative tests (rows 2-7 become Qualitative, rows 8-13 become
#    Quantitative), keeping each row's original formatting attached to its
#    content.
# ---------------------------------------------------------------------------
$wsTest = $wb.Worksheets.Item("TestDetails")

$wsTest.Range("A2").Value = "Qualitative Test -1"
$wsTest.Range("B2").Value = "Qualitative"
$wsTest.Range("A3").Value = "Qualitative Test -2"
$wsTest.Range("B3").Value = "Qualitative"
$wsTest.Range("A4").Value = "Qualitative Test -3"
$wsTest.Range("B4").Value = "Qualitative"
$wsTest.Range("A5").Value = "Qualitative Test -4"
$wsTest.Range("B5").Value = "Qualitative"
$wsTest.Range("A6").Value = "Qualitative Test -5"
$wsTest.Range("B6").Value = "Qualitative"
$wsTest.Range("A7").Value = "Qualitative Test -6"
$wsTest.Range("B7").Value = "Qualitative"

$wsTest.Range("A8").Value = "Quantitative Test-1"
$wsTest.Range("B8").Value = "Quantitative"
$wsTest.Range("A9").Value = "Quantitative Test-2"
$wsTest.Range("B9").Value = "Quantitative"
$wsTest.Range("A10").Value = "Quantitative Test-3"
$wsTest.Range("B10").Value = "Quantitative"
$wsTest.Range("A11").Value = "Quantitative Test-4"
$wsTest.Range("B11").Value = "Quantitative"
$wsTest.Range("A12").Value = "Quantitative Test-5"
$wsTest.Range("B12").Value = "Quantitative"
$wsTest.Range("A13").Value = "Quantitative Test-6"
$wsTest.Range("B13").Value = "Quantitative"

# Re-apply the original per-content font color: in the source workbook the
# "TestType" cells for Quantitative rows use a custom dark-gray font color
# (RGB 202124) while Qualitative rows use the default (theme) font color.
# Since the rows have been reordered, reapply that coloring based on the new
# row content/order: Qualitative now occupies rows 2-7 (theme font color)
# and Quantitative now occupies rows 8-13 (custom RGB font color).
$wsTest.Range("B2:B7").Font.ThemeColor = 1
$wsTest.Range("B8:B13").Font.Color = 2367776

$wb.Save()
